## T03 - OSLO - Casos de Pruebas Integracion Testify.xlsx
## "Documentacion - Gestion SQA" commit:
##   - Corrige ortografia/contenido de la portada (2da -> 3ra tanda de casos de uso,
##     "Iteracion N 3, Fase Construccion" -> "Fase Transicion")
##   - Corrige las fechas de creacion/ejecucion de los casos de prueba de transicion
##   - Deja la hoja "Portada" como hoja activa / seleccionada

$wb = $excel.ActiveWorkbook

$wsPortada = $wb.Worksheets.Item("Portada")
$wsCU06    = $wb.Worksheets.Item("CU 06 Exportar escenario y Rdo")
$wsCU08    = $wb.Worksheets.Item("CU 08 Consultar proyecto asigna")

## ---------------------------------------------------------------
## 1) Portada: correccion de textos de la caratula
## ---------------------------------------------------------------
## "2° Tanda de Casos de Uso" -> "3° Tanda de Casos de Uso" (pasaron a la fase de Transicion)
$wsPortada.Range("B9").Value = "3° Tanda de Casos de Uso"
## "Iteración N° 3, Fase Contrucción" -> "Fase Transición"
$wsPortada.Range("B12").Value = "Fase Transición"

## ---------------------------------------------------------------
## 2) CU 06 Exportar escenario y resultados: correccion de fechas
##    Fecha Creacion (E) 20-mar-2024 -> 22-mar-2025
##    Fecha Ejecucion (F) 21-mar-2024 -> 23-mar-2025
## ---------------------------------------------------------------
$cu06Rows = @(4,5,6,7,8,9,10,11,12,13,14,15,17,18,19,20)
foreach ($r in $cu06Rows) {
    $wsCU06.Range("E$r").Value = 45738
    $wsCU06.Range("F$r").Value = 45739
}

## ---------------------------------------------------------------
## 3) CU 08 Consultar proyecto asignado: correccion de fechas
##    Fecha Creacion (E) 7-nov-2024 -> 22-mar-2025
##    Fecha Ejecucion (F) 7-nov-2024 -> 23-mar-2025
## ---------------------------------------------------------------
$cu08Rows = @(4,5,6,7)
foreach ($r in $cu08Rows) {
    $wsCU08.Range("E$r").Value = 45738
    $wsCU08.Range("F$r").Value = 45739
}

## ---------------------------------------------------------------
## 4) Estado de la ventana: "Portada" pasa a ser la hoja activa
## ---------------------------------------------------------------
$wsCU06.Activate()
$wsCU06.Range("G15").Select()

$wsCU08.Activate()
$wsCU08.Range("C4").Select()

$wsPortada.Activate()
$wsPortada.Range("G10").Select()
